$d = $word.ActiveDocument

# Curly quote characters used elsewhere in this document.
$ldq = [char]8220   # “
$rdq = [char]8221   # ”

# 1. Locate the end of the final sentence in the last paragraph
#    ("...U kunt dus geen 1a, 1b of dergelijken gebruiken. ") and
#    append the new explanatory sentence there, split across two
#    runs exactly like the authored edit:
#       "Een vraag zelf moet eruitzien als volgt: “1) Vraag” of “2. Vraag"
#       "”."
$r = $d.Content
$found = $r.Find.Execute("1a, 1b of dergelijken gebruiken. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r.Collapse(0)
    $r.InsertAfter("Een vraag zelf moet eruitzien als volgt: " + $ldq + "1) Vraag" + $rdq + " of " + $ldq + "2. Vraag")
    $r.Collapse(0)
    $r.InsertAfter($rdq + ".")
    $r.Collapse(0)

    # 2. Move the "_GoBack" bookmark to sit right after the text we
    #    just typed (this is what Word does automatically whenever you
    #    edit a document - the bookmark tracks the most recent edit
    #    location). A collapsed range that lands exactly on a paragraph
    #    boundary cannot be bookmarked directly here, so a throwaway
    #    character is used to give the range real width, the bookmark
    #    is anchored around it, and then the throwaway character is
    #    removed again - the bookmark collapses back down but keeps its
    #    (correct) position.
    $r.InsertAfter("x")
    $d.Bookmarks.Add("_GoBack", $r)
    $dummy = $d.Range($r.End - 1, $r.End)
    $dummy.Delete()
}
